$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Extend the table (Table2) by one row so it covers A1:E18
# ------------------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.ListRows.Add() | Out-Null

# ------------------------------------------------------------------
# Copy the formatting of the last existing data row (17) down into
# the newly added row (18) so fills/styles (e.g. the "Medium" orange
# fill on column B, and the hyperlink style on column E) match.
# ------------------------------------------------------------------
$ws.Range("A17:E17").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Populate the new row's values (note: entered in A, E, D, B, C order
# to reproduce the exact shared-string table ordering seen in the
# target workbook: 58=Question, 59=Link, 60=Notes)
# ------------------------------------------------------------------
$ws.Range("A18").Value = "550. Game Play Analysis IV"
$ws.Range("E18").Value = "https://leetcode.com/problems/game-play-analysis-iv/solutions/3857392/100-easy-fast-clean-solution/?envType=study-plan-v2&envId=top-sql-50 "
$ws.Range("D18").Value = "Left outer join the same table on a2.player_id = a1.player_id and a2.event_date = a1.event_date + 1 where (subquery to filter) the min(event_date) of player_id in a1, grouped by player_id. The subquery is to ensure that the main query starts from the first login date for each player."
$ws.Range("B18").Value = "Medium"
$ws.Range("C18").Value = "Basic Aggregate Functions"

# ------------------------------------------------------------------
# Add the hyperlink on the new Link cell (E18)
# ------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("E18"), "https://leetcode.com/problems/game-play-analysis-iv/solutions/3857392/100-easy-fast-clean-solution/?envType=study-plan-v2&envId=top-sql-50 ") | Out-Null

# Re-apply the row's formatting so the hyperlink cell keeps the same
# visual style used by the rest of the Link column.
$ws.Range("A17:E17").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Update the active selection like the saved workbook shows
# ------------------------------------------------------------------
$ws.Range("E31").Select() | Out-Null
